# "Added manual functionality using 'shifts' (D-pad)" -----------------------
# This script edits the "Joystick Mapping" workbook (Sheet1 is the active
# joystick-mapping sheet) to describe the new D-pad "shift" based control
# scheme instead of the old position-select scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update the date stamp in the top-right corner.
# ---------------------------------------------------------------------
$ws.Range("G1").Value2 = "(Feb 18, 2019)"

# ---------------------------------------------------------------------
# 2. Operator joystick button labels (left column, rows 3-12).
# ---------------------------------------------------------------------
$ws.Range("A3").Value2 = "Deploy Ball Intake (Do NOT use)"
$ws.Range("A5").Value2 = "Expand/Contract Hatch Panel Grabber"
$ws.Range("A6").Value2 = "Shift 1"
$ws.Range("A7").Value2 = "Shift 4"
$ws.Range("A8").Value2 = "Shift 3"
$ws.Range("A9").Value2 = "Shift 2"
$ws.Range("A10").Value2 = ""
$ws.Range("A11").Value2 = "Manipulator Intake In/Out"
$ws.Range("A12").Value2 = "Kill Everything"

# ---------------------------------------------------------------------
# 3. Right-hand annotation column - remove the old per-button notes that
#    no longer apply and add the new shift explanations, merged across
#    the two rows they now span, in a smaller word-wrapped font.
# ---------------------------------------------------------------------
$ws.Range("G5").Value2 = ""
$ws.Range("G6").Value2 = ""
$ws.Range("G7").Value2 = ""
$ws.Range("G8").Value2 = ""
$ws.Range("G9").Value2 = ""
$ws.Range("G11").Value2 = ""
$ws.Range("G12").Value2 = "Kill Everything"

$ws.Range("G4").Value2 = "No Shift = Enable/Disable Ball Intake Wheels Shift 1 = Ball Intake Unfold/Fold                   Shift 2 = Shift To High/Low Gear"
$ws.Range("G4").Font.Size = 10
$ws.Range("G4").WrapText = $true
$ws.Range("G4").VerticalAlignment = -4107
$ws.Range("G4").HorizontalAlignment = 1
$ws.Range("G4:G5").Merge()

$ws.Range("G10").Value2 = "Vertical Axis:                                                         No Shift = Elevator Up/Down                                   Shift 1 = Ball Intake In/Out                                   Shift 2 = Manipulator Flip Fwd/Back"
$ws.Range("G10").Font.Size = 10
$ws.Range("G10").WrapText = $true
$ws.Range("G10").HorizontalAlignment = -4131
$ws.Range("G10").VerticalAlignment = -4107
$ws.Range("G10:G11").Merge()

# A couple of helper cells the new layout reaches into (H10, J4:J8) just to
# keep the used-range / selection highlighting lined up with the new notes.
$ws.Range("H10").Font.Bold = $false
$ws.Range("J4:J8").Font.Bold = $true

# ---------------------------------------------------------------------
# 4. Clear out the now-unused position/rocket-level reference table
#    (rows 14-18) - the whole "shift" mechanism replaces it.
# ---------------------------------------------------------------------
$ws.Range("B14:F18").Value2 = ""

# ---------------------------------------------------------------------
# 5. Driver joystick button labels (bottom half, rows 29-31).
# ---------------------------------------------------------------------
$ws.Range("G29").Value2 = "Arcade Drive Turn"
$ws.Range("A30").Value2 = "Arcade Drive Fwd/Back"
$ws.Range("A31").Value2 = "Kill Everything"
$ws.Range("G31").Value2 = "Kill Everything"

# ---------------------------------------------------------------------
# 6. Widen columns A and G to fit the longer shift descriptions and move
#    the active selection to A5, matching the saved workbook view.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.7109375
$ws.Columns.Item(7).ColumnWidth = 35.7109375
$ws.Range("A5").Select()
